$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.790254473686218
$ws.Range("B1").Value = 2.249180316925049
$ws.Range("C1").Value = 1.962173819541931
$ws.Range("D1").Value = 1.622416496276855
$ws.Range("E1").Value = 1.529477834701538
